$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Full roster table (rows 2-19) reflecting the reordering / replacement
# described by the commit (Malik Monk / Sacramento Kings dropped in favor
# of Khris Middleton, and several players reshuffled to correct rows).
$data = @(
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Wendell Carter Jr.", "C", "Orlando Magic"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Kyle Kuzma", "PF", "Washington Wizards")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
